# Update scripts with new TPM-derived statistics for the Gnai2-Cxcr2
# ligand-receptor pair table (NATMI output). Columns G:T on rows 2-9
# hold recalculated expression / specificity values produced by the
# updated analysis script; overwrite them with the refreshed numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 201.4397426666667
$ws.Range("H2").Value = 604.3192280000001
$ws.Range("I2").Value = 0.4833500233086392
$ws.Range("J2").Value = 0.4833500233086393
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.03927866666666666
$ws.Range("N2").Value = 0.117836
$ws.Range("O2").Value = 0.7432525340448212
$ws.Range("P2").Value = 0.7432525340448213
$ws.Range("Q2").Value = 7.912284505623111
$ws.Range("R2").Value = 71.21056055060801
$ws.Range("S2").Value = 0.3592511296547695
$ws.Range("T2").Value = 0.3592511296547696

$ws.Range("G3").Value = 201.4397426666667
$ws.Range("H3").Value = 604.3192280000001
$ws.Range("I3").Value = 0.4833500233086392
$ws.Range("J3").Value = 0.4833500233086393
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.01356833333333333
$ws.Range("N3").Value = 0.040705
$ws.Range("O3").Value = 0.2567474659551788
$ws.Range("P3").Value = 0.2567474659551788
$ws.Range("Q3").Value = 2.733201575082223
$ws.Range("R3").Value = 24.59881417574
$ws.Range("S3").Value = 0.1240988936538697
$ws.Range("T3").Value = 0.1240988936538697

$ws.Range("I4").Value = 0.1569674599353791
$ws.Range("J4").Value = 0.1569674599353792
$ws.Range("M4").Value = 0.03927866666666666
$ws.Range("N4").Value = 0.117836
$ws.Range("O4").Value = 0.7432525340448212
$ws.Range("P4").Value = 0.7432525340448213
$ws.Range("Q4").Value = 2.569506860953778
$ws.Range("R4").Value = 23.125561748584
$ws.Range("S4").Value = 0.1166664623595495
$ws.Range("T4").Value = 0.1166664623595495

$ws.Range("I5").Value = 0.1569674599353791
$ws.Range("J5").Value = 0.1569674599353792
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.01356833333333333
$ws.Range("N5").Value = 0.040705
$ws.Range("O5").Value = 0.2567474659551788
$ws.Range("P5").Value = 0.2567474659551788
$ws.Range("Q5").Value = 0.8876046095855556
$ws.Range("R5").Value = 7.988441486269999
$ws.Range("S5").Value = 0.04030099757582965
$ws.Range("T5").Value = 0.04030099757582965

$ws.Range("G6").Value = 60.43484133333334
$ws.Range("H6").Value = 181.304524
$ws.Range("I6").Value = 0.1450120099461104
$ws.Range("J6").Value = 0.1450120099461104
$ws.Range("M6").Value = 0.03927866666666666
$ws.Range("N6").Value = 0.117836
$ws.Range("O6").Value = 0.7432525340448212
$ws.Range("P6").Value = 0.7432525340448213
$ws.Range("Q6").Value = 2.373799987784889
$ws.Range("R6").Value = 21.364199890064
$ws.Range("S6").Value = 0.1077805438593793
$ws.Range("T6").Value = 0.1077805438593794

$ws.Range("G7").Value = 60.43484133333334
$ws.Range("H7").Value = 181.304524
$ws.Range("I7").Value = 0.1450120099461104
$ws.Range("J7").Value = 0.1450120099461104
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.01356833333333333
$ws.Range("N7").Value = 0.040705
$ws.Range("O7").Value = 0.2567474659551788
$ws.Range("P7").Value = 0.2567474659551788
$ws.Range("Q7").Value = 0.8200000721577778
$ws.Range("R7").Value = 7.38000064942
$ws.Range("S7").Value = 0.03723146608673102
$ws.Range("T7").Value = 0.03723146608673102

$ws.Range("G8").Value = 89.46554166666668
$ws.Range("H8").Value = 268.396625
$ws.Range("I8").Value = 0.2146705068098712
$ws.Range("J8").Value = 0.2146705068098712
$ws.Range("M8").Value = 0.03927866666666666
$ws.Range("N8").Value = 0.117836
$ws.Range("O8").Value = 0.7432525340448212
$ws.Range("P8").Value = 0.7432525340448213
$ws.Range("Q8").Value = 3.514087189277778
$ws.Range("R8").Value = 31.6267847035
$ws.Range("S8").Value = 0.1595543981711228
$ws.Range("T8").Value = 0.1595543981711229

$ws.Range("G9").Value = 89.46554166666668
$ws.Range("H9").Value = 268.396625
$ws.Range("I9").Value = 0.2146705068098712
$ws.Range("J9").Value = 0.2146705068098712
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.01356833333333333
$ws.Range("N9").Value = 0.040705
$ws.Range("O9").Value = 0.2567474659551788
$ws.Range("P9").Value = 0.2567474659551788
$ws.Range("Q9").Value = 1.213898291180556
$ws.Range("R9").Value = 10.925084620625
$ws.Range("S9").Value = 0.05511610863874838
$ws.Range("T9").Value = 0.05511610863874839
